# The document contains paragraphs built from three separate runs that
# together spell out "<id>p026r_1</id>" (and, in another paragraph,
# "<id>p026r_2</id>"): a Courier-New "<id>" run, an Arial run holding the
# bare id, and a Courier-New "</id>" run. The edit collapses each of these
# triples into a single run (carrying the formatting of the first/"<id>"
# run) whose text is the full "<id>...</id>" string.
#
# We find the two target paragraphs explicitly (by their exact current
# text) and re-write their content range in one shot: Word merges a
# multi-run Range.Text assignment into a single run using the formatting
# of the range's first run, exactly matching the diff.

$d = $word.ActiveDocument

$targets = @(
    @{ Old = "<id>p026r_1</id>"; New = "<id>p026r_1</id>" },
    @{ Old = "<id>p026r_2</id>"; New = "<id>p026r_2</id>" }
)

foreach ($t in $targets) {
    foreach ($p in $d.Paragraphs) {
        $pRange = $p.Range
        $pText = $pRange.Text

        if ($pText -eq ($t.Old + "`r")) {
            $start = $pRange.Start
            $end = $pRange.End - 1
            $sub = $d.Range($start, $end)

            # Force a real content change (the text is already correct,
            # so assigning it verbatim would be a no-op and the existing
            # three runs would be left untouched). Route through a
            # throwaway placeholder first, then set the final text; this
            # makes Word actually collapse the run span into one run
            # that inherits the first run's formatting.
            $sub.Text = "ZZZ_TMP_PLACEHOLDER_ZZZ"

            $pRange2 = $p.Range
            $sub2 = $d.Range($pRange2.Start, $pRange2.End - 1)
            $sub2.Text = $t.New
        }
    }
}
